# Generate Report for Handoff
# Updates the UUID-based file names, the handoff xlf file names, and the
# handoff datetimes across the Overview / zh-cn / de-de sheets, including
# the hyperlink display text shown for each updated cell.

$wb = $excel.ActiveWorkbook

$oldMd1 = "38df7ae0-6aa5-435e-a9ad-761ba9273280.md"
$newMd1 = "21061f60-c357-4a14-9b98-4f74502a4f1e.md"

$oldMd2 = "bd18b446-8748-444e-a708-4bd16e20491a.md"
$newMd2 = "11f8ab3f-d199-4602-b75d-83766375f042.md"

$oldXlfZh = "38df7ae0-6aa5-435e-a9ad-761ba9273280.3dbbd114486b22ece18115e08a0e8a8e327d07ec.zh-cn.xlf"
$newXlfZh = "21061f60-c357-4a14-9b98-4f74502a4f1e.6e8ebfcbf68aa50ba29392b1c613b238de614c9c.zh-cn.xlf"

$oldDateZh = "2016-02-15 02:59:13"
$newDateZh = "2016-02-15 03:00:23"

$oldXlfDe = "38df7ae0-6aa5-435e-a9ad-761ba9273280.3dbbd114486b22ece18115e08a0e8a8e327d07ec.de-de.xlf"
$newXlfDe = "21061f60-c357-4a14-9b98-4f74502a4f1e.6e8ebfcbf68aa50ba29392b1c613b238de614c9c.de-de.xlf"

$oldDateDe = "2016-02-15 02:59:55"
$newDateDe = "2016-02-15 03:00:36"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Hyperlinks.Item(1).TextToDisplay = $newMd1
$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Hyperlinks.Item(2).TextToDisplay = $newMd2

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd1
$wsZh.Hyperlinks.Item(1).TextToDisplay = $newMd1
$wsZh.Range("C2").Value = $newXlfZh
$wsZh.Hyperlinks.Item(2).TextToDisplay = $newXlfZh
$wsZh.Range("D2").Value = $newDateZh

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd1
$wsDe.Hyperlinks.Item(1).TextToDisplay = $newMd1
$wsDe.Range("C2").Value = $newXlfDe
$wsDe.Hyperlinks.Item(2).TextToDisplay = $newXlfDe
$wsDe.Range("D2").Value = $newDateDe

$wb.Save()
